# Swap the data of row 2 <-> row 3, and row 4 <-> row 5
# (as a result of inserting a new article and re-bucketing the time-bucket analysis).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2, $lastCol) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-Rows $ws 2 3 5
Swap-Rows $ws 4 5 5
